$wb = $excel.ActiveWorkbook

# This script applies cached-value updates (market price refresh) to the
# "Profits" sheets, mirroring a scheduled data-refresh run. Each leve/recipe
# row has columns H:N holding prices/profits pulled from market data; only
# the numeric values change (no formulas are present in these cells).

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 868
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 802
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 802
$ws.Range("M32").Value = -674
$ws.Range("N32").Value = -1454
$ws.Range("H64").Value = 3076.3635
$ws.Range("I64").Value = 3225
$ws.Range("J64").Value = 2991.4285
$ws.Range("K64").Value = 3225
$ws.Range("L64").Value = 2991.4285
$ws.Range("M64").Value = -2977
$ws.Range("N64").Value = -3487.4285
$ws.Range("H67").Value = 3076.3635
$ws.Range("I67").Value = 3225
$ws.Range("J67").Value = 2991.4285
$ws.Range("K67").Value = 3225
$ws.Range("L67").Value = 2991.4285
$ws.Range("M67").Value = -2367
$ws.Range("N67").Value = -4707.4285
$ws.Range("H70").Value = 1737.0857
$ws.Range("I70").Value = 970.75
$ws.Range("J70").Value = 2758.8667
$ws.Range("K70").Value = 2912.25
$ws.Range("L70").Value = 8276.6001
$ws.Range("M70").Value = -2642.25
$ws.Range("N70").Value = -8816.6001
$ws.Range("H73").Value = 1737.0857
$ws.Range("I73").Value = 970.75
$ws.Range("J73").Value = 2758.8667
$ws.Range("K73").Value = 2912.25
$ws.Range("L73").Value = 8276.6001
$ws.Range("M73").Value = -1976.25
$ws.Range("N73").Value = -10148.6001
$ws.Range("H93").Value = 38400
$ws.Range("J93").Value = 38400
$ws.Range("L93").Value = 38400
$ws.Range("N93").Value = -43392
$ws.Range("H96").Value = 676.17645
$ws.Range("I96").Value = 436
$ws.Range("J96").Value = 889.6667
$ws.Range("K96").Value = 1308
$ws.Range("L96").Value = 2669.0001
$ws.Range("M96").Value = 65
$ws.Range("N96").Value = -5415.0001
$ws.Range("H111").Value = 1124.1428
$ws.Range("I111").Value = 1436.5
$ws.Range("J111").Value = 889.875
$ws.Range("K111").Value = 4309.5
$ws.Range("L111").Value = 2669.625
$ws.Range("M111").Value = -1242.5
$ws.Range("N111").Value = -8803.625
$ws.Range("H116").Value = 4755.077
$ws.Range("I116").Value = 5231.231
$ws.Range("J116").Value = 4278.923
$ws.Range("K116").Value = 5231.231
$ws.Range("L116").Value = 4278.923
$ws.Range("M116").Value = -1789.231
$ws.Range("N116").Value = -11162.923
$ws.Range("H121").Value = 1107
$ws.Range("I121").Value = 557.6
$ws.Range("J121").Value = 1450.375
$ws.Range("K121").Value = 1672.8
$ws.Range("L121").Value = 4351.125
$ws.Range("M121").Value = 74.19999999999982
$ws.Range("N121").Value = -7845.125
$ws.Range("H135").Value = 794.4151000000001
$ws.Range("I135").Value = 557.86664
$ws.Range("J135").Value = 2125
$ws.Range("K135").Value = 5020.79976
$ws.Range("L135").Value = 19125
$ws.Range("M135").Value = -2485.79976
$ws.Range("N135").Value = -24195
$ws.Range("H137").Value = 3482.1086
$ws.Range("I137").Value = 3525.5527
$ws.Range("J137").Value = 3275.75
$ws.Range("K137").Value = 10576.6581
$ws.Range("L137").Value = 9827.25
$ws.Range("M137").Value = -8026.658100000001
$ws.Range("N137").Value = -14927.25
$ws.Range("H138").Value = 2417.5696
$ws.Range("I138").Value = 951.74
$ws.Range("J138").Value = 4944.8623
$ws.Range("K138").Value = 2855.22
$ws.Range("L138").Value = 14834.5869
$ws.Range("M138").Value = 2284.78
$ws.Range("N138").Value = -25114.5869
$ws.Range("H141").Value = 220583.67
$ws.Range("I141").Value = 3687.093
$ws.Range("J141").Value = 2085894.2
$ws.Range("K141").Value = 11061.279
$ws.Range("L141").Value = 6257682.6
$ws.Range("M141").Value = -5881.278999999999
$ws.Range("N141").Value = -6268042.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 985.59
$ws.Range("I32").Value = 985.30304
$ws.Range("K32").Value = 985.30304
$ws.Range("M32").Value = -698.30304
$ws.Range("H74").Value = 1090.7142
$ws.Range("I74").Value = 759.63635
$ws.Range("J74").Value = 2304.6667
$ws.Range("K74").Value = 759.63635
$ws.Range("L74").Value = 2304.6667
$ws.Range("M74").Value = 114.36365
$ws.Range("N74").Value = -4052.6667
$ws.Range("H77").Value = 1090.7142
$ws.Range("I77").Value = 759.63635
$ws.Range("J77").Value = 2304.6667
$ws.Range("K77").Value = 3798.18175
$ws.Range("L77").Value = 11523.3335
$ws.Range("M77").Value = 569.8182500000003
$ws.Range("N77").Value = -20259.3335
$ws.Range("H97").Value = 1052.5
$ws.Range("I97").Value = 1052.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1052.5
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -556.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 762.8889
$ws.Range("I94").Value = 525.1429000000001
$ws.Range("J94").Value = 1595
$ws.Range("K94").Value = 525.1429000000001
$ws.Range("L94").Value = 1595
$ws.Range("M94").Value = -74.14290000000005
$ws.Range("N94").Value = -2497
$ws.Range("H99").Value = 2646.6667
$ws.Range("I99").Value = 1454.2858
$ws.Range("J99").Value = 4316
$ws.Range("K99").Value = 1454.2858
$ws.Range("L99").Value = 4316
$ws.Range("M99").Value = 43.71419999999989
$ws.Range("N99").Value = -7312
$ws.Range("H132").Value = 20389.5
$ws.Range("J132").Value = 20389.5
$ws.Range("L132").Value = 20389.5
$ws.Range("N132").Value = -30509.5
$ws.Range("H135").Value = 30000.5
$ws.Range("J135").Value = 30000.5
$ws.Range("L135").Value = 30000.5
$ws.Range("N135").Value = -40140.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3199.3096
$ws.Range("I31").Value = 1973.037
$ws.Range("J31").Value = 5406.6
$ws.Range("K31").Value = 1973.037
$ws.Range("L31").Value = 5406.6
$ws.Range("M31").Value = -1678.037
$ws.Range("N31").Value = -5996.6
$ws.Range("H34").Value = 3199.3096
$ws.Range("I34").Value = 1973.037
$ws.Range("J34").Value = 5406.6
$ws.Range("K34").Value = 1973.037
$ws.Range("L34").Value = 5406.6
$ws.Range("M34").Value = -1771.037
$ws.Range("N34").Value = -5810.6
$ws.Range("H132").Value = 2780.6785
$ws.Range("I132").Value = 2172.0557
$ws.Range("J132").Value = 3876.2
$ws.Range("K132").Value = 6516.1671
$ws.Range("L132").Value = 11628.6
$ws.Range("M132").Value = -3986.1671
$ws.Range("N132").Value = -16688.6
$ws.Range("H134").Value = 2440.4614
$ws.Range("I134").Value = 1422.4
$ws.Range("K134").Value = 4267.200000000001
$ws.Range("M134").Value = -1732.200000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3810
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 3911.25
$ws.Range("K63").Value = 9000
$ws.Range("L63").Value = 11733.75
$ws.Range("M63").Value = -8251
$ws.Range("N63").Value = -13231.75
$ws.Range("H66").Value = 3810
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 3911.25
$ws.Range("K66").Value = 27000
$ws.Range("L66").Value = 35201.25
$ws.Range("M66").Value = -23256
$ws.Range("N66").Value = -42689.25
$ws.Range("H70").Value = 3944.25
$ws.Range("J70").Value = 6010.5
$ws.Range("L70").Value = 18031.5
$ws.Range("N70").Value = -18661.5
$ws.Range("H73").Value = 3944.25
$ws.Range("J73").Value = 6010.5
$ws.Range("L73").Value = 18031.5
$ws.Range("N73").Value = -20215.5
$ws.Range("H87").Value = 7643.4546
$ws.Range("I87").Value = 2759.6365
$ws.Range("J87").Value = 12527.272
$ws.Range("K87").Value = 8278.9095
$ws.Range("L87").Value = 37581.81600000001
$ws.Range("M87").Value = -7030.9095
$ws.Range("N87").Value = -40077.81600000001
$ws.Range("H90").Value = 7643.4546
$ws.Range("I90").Value = 2759.6365
$ws.Range("J90").Value = 12527.272
$ws.Range("K90").Value = 24836.7285
$ws.Range("L90").Value = 112745.448
$ws.Range("M90").Value = -18596.7285
$ws.Range("N90").Value = -125225.448
$ws.Range("H106").Value = 3999.889
$ws.Range("J106").Value = 3999.889
$ws.Range("L106").Value = 11999.667
$ws.Range("N106").Value = -13891.667

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 31384.615
$ws.Range("J64").Value = 31384.615
$ws.Range("L64").Value = 31384.615
$ws.Range("N64").Value = -31880.615
$ws.Range("H67").Value = 31384.615
$ws.Range("J67").Value = 31384.615
$ws.Range("L67").Value = 31384.615
$ws.Range("N67").Value = -33100.61500000001
$ws.Range("H80").Value = 2566.0715
$ws.Range("I80").Value = 2672
$ws.Range("J80").Value = 2301.25
$ws.Range("K80").Value = 2672
$ws.Range("L80").Value = 2301.25
$ws.Range("M80").Value = -1674
$ws.Range("N80").Value = -4297.25
$ws.Range("H83").Value = 2566.0715
$ws.Range("I83").Value = 2672
$ws.Range("J83").Value = 2301.25
$ws.Range("K83").Value = 13360
$ws.Range("L83").Value = 11506.25
$ws.Range("M83").Value = -8368
$ws.Range("N83").Value = -21490.25
$ws.Range("H113").Value = 2256.3157
$ws.Range("I113").Value = 1978
$ws.Range("J113").Value = 3300
$ws.Range("K113").Value = 1978
$ws.Range("L113").Value = 3300
$ws.Range("M113").Value = 192
$ws.Range("N113").Value = -7640

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3560
$ws.Range("I40").Value = 1900
$ws.Range("J40").Value = 4666.6665
$ws.Range("K40").Value = 1900
$ws.Range("L40").Value = 4666.6665
$ws.Range("M40").Value = -1764
$ws.Range("N40").Value = -4938.6665
$ws.Range("H93").Value = 4000
$ws.Range("I93").Value = 3520
$ws.Range("J93").Value = 6400
$ws.Range("K93").Value = 3520
$ws.Range("L93").Value = 6400
$ws.Range("M93").Value = -2272
$ws.Range("N93").Value = -8896
$ws.Range("H132").Value = 1886.8864
$ws.Range("I132").Value = 985.125
$ws.Range("K132").Value = 2955.375
$ws.Range("M132").Value = -425.375
$ws.Range("H133").Value = 42000
$ws.Range("J133").Value = 42000
$ws.Range("L133").Value = 42000
$ws.Range("N133").Value = -47060
$ws.Range("H136").Value = 1553.8422
$ws.Range("I136").Value = 1072.8368
$ws.Range("K136").Value = 3218.5104
$ws.Range("M136").Value = -668.5104000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 34500
$ws.Range("J69").Value = 34500
$ws.Range("L69").Value = 34500
$ws.Range("N69").Value = -35998
$ws.Range("H72").Value = 34500
$ws.Range("J72").Value = 34500
$ws.Range("L72").Value = 103500
$ws.Range("N72").Value = -110988
$ws.Range("H92").Value = 39833.332
$ws.Range("J92").Value = 39833.332
$ws.Range("L92").Value = 39833.332
$ws.Range("N92").Value = -44825.332
$ws.Range("H125").Value = 24000
$ws.Range("J125").Value = 24000
$ws.Range("L125").Value = 24000
$ws.Range("N125").Value = -33840
$ws.Range("H132").Value = 3663.898
$ws.Range("I132").Value = 1254.125
$ws.Range("J132").Value = 8199.941000000001
$ws.Range("K132").Value = 3762.375
$ws.Range("L132").Value = 24599.823
$ws.Range("M132").Value = -1232.375
$ws.Range("N132").Value = -29659.823
